$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at N (14th column), shifting the "目前利率" / "調整後利率"
# columns one to the right.
$ws.Columns("N:N").Insert()

# Set the header text for the newly inserted column.
$ws.Range("N1").Value = "利率種類"

# Resize columns M and N to match the new layout (values taken from the target
# workbook; the COM layer quantizes ColumnWidth to a 1/7 character grid, so the
# closest representable width is used).
$ws.Columns("M:M").ColumnWidth = 10.428571428571429
$ws.Columns("N:N").ColumnWidth = 9.428571428571429

# Update the _FilterDatabase defined name so it covers the new last column (P).
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "=正常件!`$A`$1:`$P`$1"
    }
}

# Move the active selection to the newly inserted header cell.
[void]$ws.Range("N1").Select()
